$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - "Save", same style as the other header cells (copy from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Data values for the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
